$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.639979362487793
$ws.Range("B1").Value = 4.312362194061279
$ws.Range("C1").Value = 2.635947227478027
$ws.Range("D1").Value = 2.279310941696167
$ws.Range("E1").Value = 1.775561451911926
